# Add a new "J7" jornada (matchday) column to the partidos sheet,
# continuing the existing ACB24_104xxx match id sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "J7"
$ws.Range("H2").Value = "ACB24_104106"
$ws.Range("H3").Value = "ACB24_104107"
$ws.Range("H4").Value = "ACB24_104108"
$ws.Range("H5").Value = "ACB24_104109"
$ws.Range("H6").Value = "ACB24_104110"
$ws.Range("H7").Value = "ACB24_104111"
$ws.Range("H8").Value = "ACB24_104112"
$ws.Range("H9").Value = "ACB24_104113"
$ws.Range("H10").Value = "ACB24_104114"

$ws.Columns("H").AutoFit() | Out-Null

$ws.Range("J8").Select() | Out-Null
